$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-25 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-26 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("51+3=54", $true, $false, $false, $false, $false, $true, 1, $false, "92-37=55", 2) | Out-Null
$d.Content.Find.Execute("15-13=2", $true, $false, $false, $false, $false, $true, 1, $false, "79+4=83", 2) | Out-Null
$d.Content.Find.Execute("50-22=28", $true, $false, $false, $false, $false, $true, 1, $false, "50-37=13", 2) | Out-Null
$d.Content.Find.Execute("73-58=15", $true, $false, $false, $false, $false, $true, 1, $false, "63-60=3", 2) | Out-Null
$d.Content.Find.Execute("60+11=71", $true, $false, $false, $false, $false, $true, 1, $false, "83-47=36", 2) | Out-Null
$d.Content.Find.Execute("10+62=72", $true, $false, $false, $false, $false, $true, 1, $false, "87-83=4", 2) | Out-Null
$d.Content.Find.Execute("58-23=35", $true, $false, $false, $false, $false, $true, 1, $false, "89-63=26", 2) | Out-Null
$d.Content.Find.Execute("40+20=60", $true, $false, $false, $false, $false, $true, 1, $false, "66-1=65", 2) | Out-Null
$d.Content.Find.Execute("13+67=80", $true, $false, $false, $false, $false, $true, 1, $false, "38-17=21", 2) | Out-Null
$d.Content.Find.Execute("74-15=59", $true, $false, $false, $false, $false, $true, 1, $false, "68-29=39", 2) | Out-Null
$d.Content.Find.Execute("42+38=80", $true, $false, $false, $false, $false, $true, 1, $false, "43+21=64", 2) | Out-Null
$d.Content.Find.Execute("72-15=57", $true, $false, $false, $false, $false, $true, 1, $false, "99-34=65", 2) | Out-Null
$d.Content.Find.Execute("67-10=57", $true, $false, $false, $false, $false, $true, 1, $false, "30+60=90", 2) | Out-Null
$d.Content.Find.Execute("10+11=21", $true, $false, $false, $false, $false, $true, 1, $false, "3+23=26", 2) | Out-Null
$d.Content.Find.Execute("87-65=22", $true, $false, $false, $false, $false, $true, 1, $false, "19+4=23", 2) | Out-Null
$d.Content.Find.Execute("26-22=4", $true, $false, $false, $false, $false, $true, 1, $false, "13+20=33", 2) | Out-Null
$d.Content.Find.Execute("5+37=42", $true, $false, $false, $false, $false, $true, 1, $false, "38-25=13", 2) | Out-Null
$d.Content.Find.Execute("19-11=8", $true, $false, $false, $false, $false, $true, 1, $false, "56+33=89", 2) | Out-Null
$d.Content.Find.Execute("67+18=85", $true, $false, $false, $false, $false, $true, 1, $false, "85-21=64", 2) | Out-Null
$d.Content.Find.Execute("10+4=14", $true, $false, $false, $false, $false, $true, 1, $false, "87-20=67", 2) | Out-Null
$d.Content.Find.Execute("89-15=74", $true, $false, $false, $false, $false, $true, 1, $false, "98-76=22", 2) | Out-Null
$d.Content.Find.Execute("29+55=84", $true, $false, $false, $false, $false, $true, 1, $false, "82+17=99", 2) | Out-Null
$d.Content.Find.Execute("36+11=47", $true, $false, $false, $false, $false, $true, 1, $false, "6+43=49", 2) | Out-Null
$d.Content.Find.Execute("55+38=93", $true, $false, $false, $false, $false, $true, 1, $false, "23+45=68", 2) | Out-Null
$d.Content.Find.Execute("94-67=27", $true, $false, $false, $false, $false, $true, 1, $false, "55-37=18", 2) | Out-Null
$d.Content.Find.Execute("1+56=57", $true, $false, $false, $false, $false, $true, 1, $false, "96-55=41", 2) | Out-Null
$d.Content.Find.Execute("70-5=65", $true, $false, $false, $false, $false, $true, 1, $false, "49+8=57", 2) | Out-Null
$d.Content.Find.Execute("26+60=86", $true, $false, $false, $false, $false, $true, 1, $false, "54-26=28", 2) | Out-Null
$d.Content.Find.Execute("70+8=78", $true, $false, $false, $false, $false, $true, 1, $false, "49+34=83", 2) | Out-Null
$d.Content.Find.Execute("89+0=89", $true, $false, $false, $false, $false, $true, 1, $false, "57-42=15", 2) | Out-Null
$d.Content.Find.Execute("4+20=24", $true, $false, $false, $false, $false, $true, 1, $false, "85-49=36", 2) | Out-Null
$d.Content.Find.Execute("86-59=27", $true, $false, $false, $false, $false, $true, 1, $false, "20-4=16", 2) | Out-Null
$d.Content.Find.Execute("98-78=20", $true, $false, $false, $false, $false, $true, 1, $false, "74+3=77", 2) | Out-Null
$d.Content.Find.Execute("86-64=22", $true, $false, $false, $false, $false, $true, 1, $false, "2+95=97", 2) | Out-Null
$d.Content.Find.Execute("39+5=44", $true, $false, $false, $false, $false, $true, 1, $false, "31+12=43", 2) | Out-Null
$d.Content.Find.Execute("55-43=12", $true, $false, $false, $false, $false, $true, 1, $false, "51-6=45", 2) | Out-Null
$d.Content.Find.Execute("67+13=80", $true, $false, $false, $false, $false, $true, 1, $false, "35-28=7", 2) | Out-Null
$d.Content.Find.Execute("74-34=40", $true, $false, $false, $false, $false, $true, 1, $false, "53+1=54", 2) | Out-Null
$d.Content.Find.Execute("1+60=61", $true, $false, $false, $false, $false, $true, 1, $false, "47+40=87", 2) | Out-Null
$d.Content.Find.Execute("77-9=68", $true, $false, $false, $false, $false, $true, 1, $false, "62-32=30", 2) | Out-Null
$d.Content.Find.Execute("30-5=25", $true, $false, $false, $false, $false, $true, 1, $false, "67+0=67", 2) | Out-Null
$d.Content.Find.Execute("50+42=92", $true, $false, $false, $false, $false, $true, 1, $false, "2+63=65", 2) | Out-Null
$d.Content.Find.Execute("82-14=68", $true, $false, $false, $false, $false, $true, 1, $false, "64-6=58", 2) | Out-Null
$d.Content.Find.Execute("66-20=46", $true, $false, $false, $false, $false, $true, 1, $false, "71-41=30", 2) | Out-Null
$d.Content.Find.Execute("59+19=78", $true, $false, $false, $false, $false, $true, 1, $false, "65-27=38", 2) | Out-Null
$d.Content.Find.Execute("69-51=18", $true, $false, $false, $false, $false, $true, 1, $false, "25+53=78", 2) | Out-Null
$d.Content.Find.Execute("96-81=15", $true, $false, $false, $false, $false, $true, 1, $false, "4+26=30", 2) | Out-Null
$d.Content.Find.Execute("24+0=24", $true, $false, $false, $false, $false, $true, 1, $false, "36-34=2", 2) | Out-Null
$d.Content.Find.Execute("90-80=10", $true, $false, $false, $false, $false, $true, 1, $false, "26+3=29", 2) | Out-Null
$d.Content.Find.Execute("72-43=29", $true, $false, $false, $false, $false, $true, 1, $false, "13+16=29", 2) | Out-Null
$d.Content.Find.Execute("89-46=43", $true, $false, $false, $false, $false, $true, 1, $false, "90-61=29", 2) | Out-Null
$d.Content.Find.Execute("22+22=44", $true, $false, $false, $false, $false, $true, 1, $false, "48-3=45", 2) | Out-Null
$d.Content.Find.Execute("31+39=70", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=25", 2) | Out-Null
$d.Content.Find.Execute("45-41=4", $true, $false, $false, $false, $false, $true, 1, $false, "67-24=43", 2) | Out-Null
$d.Content.Find.Execute("18+33=51", $true, $false, $false, $false, $false, $true, 1, $false, "32-20=12", 2) | Out-Null
$d.Content.Find.Execute("43+22=65", $true, $false, $false, $false, $false, $true, 1, $false, "96-73=23", 2) | Out-Null
$d.Content.Find.Execute("45-7=38", $true, $false, $false, $false, $false, $true, 1, $false, "14+47=61", 2) | Out-Null
$d.Content.Find.Execute("83-81=2", $true, $false, $false, $false, $false, $true, 1, $false, "70-31=39", 2) | Out-Null
$d.Content.Find.Execute("30+18=48", $true, $false, $false, $false, $false, $true, 1, $false, "70-60=10", 2) | Out-Null
$d.Content.Find.Execute("10+50=60", $true, $false, $false, $false, $false, $true, 1, $false, "89-63=26", 2) | Out-Null
$d.Content.Find.Execute("60-12=48", $true, $false, $false, $false, $false, $true, 1, $false, "62-50=12", 2) | Out-Null
$d.Content.Find.Execute("64-61=3", $true, $false, $false, $false, $false, $true, 1, $false, "69-28=41", 2) | Out-Null
$d.Content.Find.Execute("82-31=51", $true, $false, $false, $false, $false, $true, 1, $false, "41-32=9", 2) | Out-Null
$d.Content.Find.Execute("47+31=78", $true, $false, $false, $false, $false, $true, 1, $false, "57-31=26", 2) | Out-Null
$d.Content.Find.Execute("48+46=94", $true, $false, $false, $false, $false, $true, 1, $false, "17+45=62", 2) | Out-Null
$d.Content.Find.Execute("67-55=12", $true, $false, $false, $false, $false, $true, 1, $false, "80-52=28", 2) | Out-Null
$d.Content.Find.Execute("15+16=31", $true, $false, $false, $false, $false, $true, 1, $false, "12-6=6", 2) | Out-Null
$d.Content.Find.Execute("59+6=65", $true, $false, $false, $false, $false, $true, 1, $false, "98-40=58", 2) | Out-Null
$d.Content.Find.Execute("17-12=5", $true, $false, $false, $false, $false, $true, 1, $false, "38+5=43", 2) | Out-Null
$d.Content.Find.Execute("9+28=37", $true, $false, $false, $false, $false, $true, 1, $false, "7+64=71", 2) | Out-Null
$d.Content.Find.Execute("84-29=55", $true, $false, $false, $false, $false, $true, 1, $false, "49-17=32", 2) | Out-Null
$d.Content.Find.Execute("19-7=12", $true, $false, $false, $false, $false, $true, 1, $false, "73-65=8", 2) | Out-Null
$d.Content.Find.Execute("62+16=78", $true, $false, $false, $false, $false, $true, 1, $false, "72-69=3", 2) | Out-Null
$d.Content.Find.Execute("25+51=76", $true, $false, $false, $false, $false, $true, 1, $false, "1+76=77", 2) | Out-Null
$d.Content.Find.Execute("51-32=19", $true, $false, $false, $false, $false, $true, 1, $false, "17+11=28", 2) | Out-Null
$d.Content.Find.Execute("44-1=43", $true, $false, $false, $false, $false, $true, 1, $false, "21+19=40", 2) | Out-Null
$d.Content.Find.Execute("22+72=94", $true, $false, $false, $false, $false, $true, 1, $false, "0+91=91", 2) | Out-Null
$d.Content.Find.Execute("22+27=49", $true, $false, $false, $false, $false, $true, 1, $false, "70+10=80", 2) | Out-Null
$d.Content.Find.Execute("71+2=73", $true, $false, $false, $false, $false, $true, 1, $false, "82-76=6", 2) | Out-Null
$d.Content.Find.Execute("81+2=83", $true, $false, $false, $false, $false, $true, 1, $false, "11+46=57", 2) | Out-Null
$d.Content.Find.Execute("83-78=5", $true, $false, $false, $false, $false, $true, 1, $false, "32+41=73", 2) | Out-Null
$d.Content.Find.Execute("83-33=50", $true, $false, $false, $false, $false, $true, 1, $false, "10+52=62", 2) | Out-Null
$d.Content.Find.Execute("94-30=64", $true, $false, $false, $false, $false, $true, 1, $false, "14+47=61", 2) | Out-Null
$d.Content.Find.Execute("58+39=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-22=67", 2) | Out-Null
$d.Content.Find.Execute("19+40=59", $true, $false, $false, $false, $false, $true, 1, $false, "38+40=78", 2) | Out-Null
$d.Content.Find.Execute("34+55=89", $true, $false, $false, $false, $false, $true, 1, $false, "79-69=10", 2) | Out-Null
$d.Content.Find.Execute("51-5=46", $true, $false, $false, $false, $false, $true, 1, $false, "85-58=27", 2) | Out-Null
$d.Content.Find.Execute("15+7=22", $true, $false, $false, $false, $false, $true, 1, $false, "99-77=22", 2) | Out-Null
$d.Content.Find.Execute("57-26=31", $true, $false, $false, $false, $false, $true, 1, $false, "86-39=47", 2) | Out-Null
$d.Content.Find.Execute("72-40=32", $true, $false, $false, $false, $false, $true, 1, $false, "1+33=34", 2) | Out-Null
$d.Content.Find.Execute("0+67=67", $true, $false, $false, $false, $false, $true, 1, $false, "40-16=24", 2) | Out-Null
$d.Content.Find.Execute("87-39=48", $true, $false, $false, $false, $false, $true, 1, $false, "29-20=9", 2) | Out-Null
$d.Content.Find.Execute("75-30=45", $true, $false, $false, $false, $false, $true, 1, $false, "0+0=0", 2) | Out-Null
$d.Content.Find.Execute("20+23=43", $true, $false, $false, $false, $false, $true, 1, $false, "24+48=72", 2) | Out-Null
$d.Content.Find.Execute("65+11=76", $true, $false, $false, $false, $false, $true, 1, $false, "19+33=52", 2) | Out-Null
$d.Content.Find.Execute("3+15=18", $true, $false, $false, $false, $false, $true, 1, $false, "57-3=54", 2) | Out-Null
$d.Content.Find.Execute("85-80=5", $true, $false, $false, $false, $false, $true, 1, $false, "6+3=9", 2) | Out-Null
$d.Content.Find.Execute("83-82=1", $true, $false, $false, $false, $false, $true, 1, $false, "71-51=20", 2) | Out-Null
$d.Content.Find.Execute("78-74=4", $true, $false, $false, $false, $false, $true, 1, $false, "77+12=89", 2) | Out-Null
$d.Content.Find.Execute("91-32=59", $true, $false, $false, $false, $false, $true, 1, $false, "99-72=27", 2) | Out-Null

Write-Output "Replacements applied: 101"
